$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (B and C) ---
# The engine's COM ColumnWidth setter stores width+5/6 into the OOXML `width`
# attribute, so subtract 5/6 from the desired stored width to land exactly
# on the target values (18 and 46).
$ws.Columns.Item(2).ColumnWidth = 18 - 5/6
$ws.Columns.Item(3).ColumnWidth = 46 - 5/6

# --- New row 20: "Course Description" ---
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = "Course Description"
$ws.Range("C20").Value = "Place background images against the courses in course description screen."
$ws.Range("C20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 30

# --- New row 21: "Site - Color Change" ---
$ws.Range("A21").Value = 12
$ws.Range("B21").Value = "Site - Color Change"
$ws.Range("C21").Value = "Now we have changed the logo color.`nSo we need to change the site color also."
$ws.Range("C21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 30

# --- Selection / view state ---
$ws.Range("C20").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
